$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data through row 97 (A1:B97). Append a new
# row 98 with the next day's gold-price entry, reusing the same
# formatting (bordered cell for the date, bordered + wrap-text cell
# for the price description) as the preceding row.
$ws.Range("A97:B97").Copy()
$ws.Range("A98:B98").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(98, 1).Value = "23-12-2025"
$ws.Cells.Item(98, 2).Value = "The price of gold in India today is ₹13,855 per gram for 24 karat gold, ₹12,700 per gram for 22 karat gold and ₹10,391 per gram for 18 karat gold (also called 999 gold)."
